$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.199.56"
$ws.Range("E2").Value = "'  -4.09%  "

$ws.Range("D3").Value = "'1.658.35"
$ws.Range("E3").Value = "'  -2.83%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "'  +0.12%  "

$ws.Range("D5").Value = "'218.02"
$ws.Range("E5").Value = "'  -2.69%  "

$ws.Range("D6").Value = "'0.5151"
$ws.Range("E6").Value = "'  -2.94%  "

$ws.Range("E7").Value = "'  +0.08%  "

$ws.Range("D8").Value = "'0.2584"
$ws.Range("E8").Value = "'  -2.91%  "

$ws.Range("D9").Value = "'0.06442"
$ws.Range("E9").Value = "'  -2.16%  "

$ws.Range("D10").Value = "'19.97"
$ws.Range("E10").Value = "'  -3.78%  "

$ws.Range("E11").Value = "'  +2.30%  "

$ws.Range("D12").Value = "'1.658.74"
$ws.Range("E12").Value = "'  -2.78%  "

$ws.Range("D13").Value = "'4.297"
$ws.Range("E13").Value = "'  -4.79%  "

$ws.Range("D14").Value = "'1.886.96"
$ws.Range("E14").Value = "'  -2.76%  "

$ws.Range("D15").Value = "'0.5548"
$ws.Range("E15").Value = "'  -3.95%  "

$ws.Range("D16").Value = "'0.0₅8075"
$ws.Range("E16").Value = "'  -0.80%  "

$ws.Range("D17").Value = "'64.32"
$ws.Range("E17").Value = "'  -4.88%  "

$ws.Range("D18").Value = "'26.225.08"
$ws.Range("E18").Value = "'  -4.00%  "

$ws.Range("D19").Value = "'212.37"
$ws.Range("E19").Value = "'  -1.49%  "

$ws.Range("E20").Value = "'  -0.01%  "

$ws.Range("D21").Value = "'4.427"
$ws.Range("E21").Value = "'  -4.11%  "

$ws.Range("E22").Value = "'  -3.20%  "

$ws.Range("D23").Value = "'5.958"
$ws.Range("E23").Value = "'  -0.04%  "

$ws.Range("E24").Value = "'  +0.02%  "

$ws.Range("D25").Value = "'143.84"
$ws.Range("E25").Value = "'  -0.43%  "

$ws.Range("D26").Value = "'1.754"
$ws.Range("E26").Value = "'  +2.96%  "

$ws.Range("D27").Value = "'0.1165"
$ws.Range("E27").Value = "'  -3.12%  "

$ws.Range("D28").Value = "'6.969"
$ws.Range("E28").Value = "'  -3.40%  "

$ws.Range("E29").Value = "'  -1.95%  "

$ws.Range("D30").Value = "'0.05250"
$ws.Range("E30").Value = "'  -2.31%  "

$ws.Range("D31").Value = "'1.254"
$ws.Range("E31").Value = "'  -2.64%  "

$ws.Range("D32").Value = "'3.368"
$ws.Range("E32").Value = "'  -2.85%  "

$ws.Range("D33").Value = "'3.223"
$ws.Range("E33").Value = "'  -5.26%  "

$ws.Range("D34").Value = "'1.570"
$ws.Range("E34").Value = "'  -4.61%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'0.9302"
$ws.Range("E36").Value = "'  -1.67%  "

$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.369"
$ws.Range("E37").Value = "'  -2.03%  "

$ws.Range("D38").Value = "'1.166.51"
$ws.Range("E38").Value = "'  +12.10%  "

$ws.Range("D39").Value = "'0.5690"
$ws.Range("E39").Value = "'  -1.86%  "

$ws.Range("D40").Value = "'0.01593"
$ws.Range("E40").Value = "'  -1.94%  "

$ws.Range("D41").Value = "'0.8483"
$ws.Range("E41").Value = "'  +0.84%  "

$ws.Range("D42").Value = "'1.004"
$ws.Range("E42").Value = "'  +0.00%  "

$ws.Range("D43").Value = "'5.681"
$ws.Range("E43").Value = "'  -1.52%  "

$ws.Range("D44").Value = "'100.39"
$ws.Range("E44").Value = "'  -0.58%  "

$ws.Range("D45").Value = "'1.796.92"
$ws.Range("E45").Value = "'  -2.79%  "

$ws.Range("E46").Value = "'  -2.91%  "

$ws.Range("D47").Value = "'0.4535"
$ws.Range("E47").Value = "'  +0.37%  "

$ws.Range("D48").Value = "'55.98"
$ws.Range("E48").Value = "'  -3.05%  "

$ws.Range("E49").Value = "'  +0.18%  "

$ws.Range("D50").Value = "'7.849"
$ws.Range("E50").Value = "'  -2.40%  "

$ws.Range("D51").Value = "'0.05056"
$ws.Range("E51").Value = "'  -3.32%  "
